$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, [string]$text) {
    $helper = $ws.Range("Z1000")
    $escaped = $text -replace '"', '""'
    $helper.Formula = '="' + $escaped + '"'
    $helper.Copy()
    $range.PasteSpecial(-4163)  # xlPasteValues
    $helper.Clear()
}

# Header updates
$ws.Range("B1").Value = "2014-08-01-2014-08-15"
$ws.Range("B3").Value = "Semi-monthly"

# New row 5
$ws.Range("A5").Value = "Web Developer1"
$ws.Range("B5").Value = 6
$ws.Range("C5").Value = "Jen Dee  Dela Cruz"
$ws.Range("D5").Value = 24
Set-TextValue $ws.Range("E5") "12,500.00"
Set-TextValue $ws.Range("F5") "1,153.85"
$ws.Range("G5").Value = "S0"
$ws.Range("H5").Value = 0
Set-TextValue $ws.Range("I5") "12,500.00"
$ws.Range("J5").Value = 290.65
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 50
$ws.Range("M5").Value = 12692.30769230769
$ws.Range("N5").Value = 0
$ws.Range("O5").Value = 0
$ws.Range("P5").Value = 0
$ws.Range("Q5").Value = -532.96
